$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The prompt text in E5 dropped its leading "item x/" — the template now
# starts directly with "{question}/{answer}.".
$ws.Range("E5").Value = "Provide a list of 20 most related best questions with answers, in this format: `n`n{question}/{answer}.`n`nFinal output are in the following format:`n    - item 1`n    - item 2`n    - item 3"

# The view was scrolled down a bit (so row 5 sits at the top of the
# visible window) and the selection moved from E6 to F5.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F5").Select()
